$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.008999999999999999
$ws.Range("B1").Value = 0.103
$ws.Range("C1").Value = -0.328
$ws.Range("D1").Value = 86
$ws.Range("E1").Value = -6
$ws.Range("F1").Value = -161

$ws.Range("A2").Value = -0.329
$ws.Range("B2").Value = -0.235
$ws.Range("C2").Value = 0.11
$ws.Range("D2").Value = -74
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 60

$ws.Range("A3").Value = 0.111
$ws.Range("B3").Value = -0.131
$ws.Range("C3").Value = 0.464
$ws.Range("D3").Value = 137
$ws.Range("E3").Value = -78
$ws.Range("F3").Value = -61

$ws.Range("A4").Value = -0.015
$ws.Range("B4").Value = -0.379
$ws.Range("C4").Value = -0.145
$ws.Range("D4").Value = -50
$ws.Range("E4").Value = -59
$ws.Range("F4").Value = 59

$ws.Range("A5").Value = -0.135
$ws.Range("B5").Value = 0.189
$ws.Range("C5").Value = -0.044
$ws.Range("D5").Value = 154
$ws.Range("E5").Value = -4
$ws.Range("F5").Value = 11

$ws.Range("A6").Value = -0.3
$ws.Range("B6").Value = 0.109
$ws.Range("C6").Value = 0.457
$ws.Range("D6").Value = -130
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 61

$ws.Range("A7").Value = 0.04
$ws.Range("B7").Value = 0.244
$ws.Range("C7").Value = 0.149
$ws.Range("D7").Value = -94
$ws.Range("E7").Value = -5
$ws.Range("F7").Value = 91

$ws.Range("A8").Value = -0.336
$ws.Range("B8").Value = 0.019
$ws.Range("C8").Value = 0.399
$ws.Range("D8").Value = -112
$ws.Range("E8").Value = -10
$ws.Range("F8").Value = -134

$ws.Range("A9").Value = -0.057
$ws.Range("B9").Value = -0.342
$ws.Range("C9").Value = 0.254
$ws.Range("D9").Value = -121
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 10

$ws.Range("A10").Value = -0.01
$ws.Range("B10").Value = 0.11
$ws.Range("C10").Value = 0.421
$ws.Range("D10").Value = -36
$ws.Range("E10").Value = -44
$ws.Range("F10").Value = -37
